$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values that changed
$ws.Range("A13").Value = 0
$ws.Range("A15").Value = 1
$ws.Range("A16").Value = 0
$ws.Range("A20").Value = 1
$ws.Range("A21").Value = 0
$ws.Range("A25").Value = 1

# Append new rows 46-48 with value 0
$ws.Range("A46").Value = 0
$ws.Range("A47").Value = 0
$ws.Range("A48").Value = 0

# Update selection/active cell and scroll position
$ws.Range("A25").Select()
